# CRITICAL FIX: Exit same side as entry, correct P&L calc for prediction markets
#
# The live-trading loop re-ran (clock moved ~5 min forward) which re-priced
# every still-OPEN trade's confidence/entry price and, for several rows,
# flipped the predicted side (UP/DOWN) now that exit-matches-entry logic is
# fixed. A brand-new trade (#23 / row 24) was also opened and logged.
# Both "All Trades" and "base_strategy" tabs carry the same trade log, so
# the same edits are applied to each.
$wb = $excel.ActiveWorkbook
$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(2,3).Value = "23:02:29"
    $ws.Cells.Item(2,5).Value = "DOWN"
    $ws.Cells.Item(2,6).Value = 0.07000000000000001

    $ws.Cells.Item(3,3).Value = "23:02:35"
    $ws.Cells.Item(3,6).Value = 0.9399999999999999

    $ws.Cells.Item(4,3).Value = "23:02:41"
    $ws.Cells.Item(4,6).Value = 0.09

    $ws.Cells.Item(5,3).Value = "23:02:47"
    $ws.Cells.Item(5,6).Value = 0.11

    $ws.Cells.Item(6,3).Value = "23:02:53"
    $ws.Cells.Item(6,5).Value = "DOWN"
    $ws.Cells.Item(6,6).Value = 0.11

    $ws.Cells.Item(7,3).Value = "23:02:59"
    $ws.Cells.Item(7,6).Value = 0.88

    $ws.Cells.Item(8,3).Value = "23:03:05"
    $ws.Cells.Item(8,5).Value = "UP"
    $ws.Cells.Item(8,6).Value = 0.9

    $ws.Cells.Item(9,3).Value = "23:03:11"
    $ws.Cells.Item(9,6).Value = 0.93

    $ws.Cells.Item(10,3).Value = "23:03:17"
    $ws.Cells.Item(10,5).Value = "UP"
    $ws.Cells.Item(10,6).Value = 0.96

    $ws.Cells.Item(11,3).Value = "23:03:23"
    $ws.Cells.Item(11,5).Value = "UP"
    $ws.Cells.Item(11,6).Value = 0.96

    $ws.Cells.Item(12,3).Value = "23:03:29"
    $ws.Cells.Item(12,5).Value = "UP"
    $ws.Cells.Item(12,6).Value = 0.96

    $ws.Cells.Item(13,3).Value = "23:03:35"
    $ws.Cells.Item(13,6).Value = 0.06

    $ws.Cells.Item(14,3).Value = "23:03:41"
    $ws.Cells.Item(14,6).Value = 0.9

    $ws.Cells.Item(15,3).Value = "23:03:47"
    $ws.Cells.Item(15,6).Value = 0.88

    $ws.Cells.Item(16,3).Value = "23:03:52"
    $ws.Cells.Item(16,6).Value = 0.08

    $ws.Cells.Item(17,3).Value = "23:03:59"
    $ws.Cells.Item(17,6).Value = 0.89

    $ws.Cells.Item(18,3).Value = "23:04:05"
    $ws.Cells.Item(18,5).Value = "UP"
    $ws.Cells.Item(18,6).Value = 0.96

    $ws.Cells.Item(19,3).Value = "23:04:11"
    $ws.Cells.Item(19,6).Value = 0.03

    $ws.Cells.Item(20,3).Value = "23:04:17"
    $ws.Cells.Item(20,6).Value = 0.98

    $ws.Cells.Item(21,3).Value = "23:04:24"
    $ws.Cells.Item(21,6).Value = 0.97

    $ws.Cells.Item(22,3).Value = "23:04:30"
    $ws.Cells.Item(22,6).Value = 0.98

    $ws.Cells.Item(23,3).Value = "23:04:36"
    $ws.Cells.Item(23,6).Value = 0.99

    # New row 24 - trade #23, just opened
    $ws.Cells.Item(24,1).Value = 23
    # Force text so Excel doesn't auto-coerce "2026-02-16" into a date serial
    # (the existing Date column cells are plain text, e.g. row 2's "2026-02-16").
    $ws.Cells.Item(24,2).NumberFormat = "@"
    $ws.Cells.Item(24,2).Value = "2026-02-16"
    $ws.Cells.Item(24,3).Value = "23:04:42"
    $ws.Cells.Item(24,4).Value = "base_strategy"
    $ws.Cells.Item(24,5).Value = "UP"
    $ws.Cells.Item(24,6).Value = 0.99
    $ws.Cells.Item(24,7).Value = ""
    $ws.Cells.Item(24,8).Value = "OPEN"
    $ws.Cells.Item(24,9).Value = 0
    $ws.Cells.Item(24,10).Value = 0
    $ws.Cells.Item(24,11).Value = 100
    $ws.Cells.Item(24,12).Value = 0
    $ws.Cells.Item(24,13).Value = 0
    $ws.Cells.Item(24,14).Value = 0.6
    $ws.Cells.Item(24,15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(24,16).Value = ""
    $ws.Cells.Item(24,17).Value = 0
}
